# "some testing set up"
#
# The authoritative diff for this commit only touches relationship-id
# values (p:sldMasterId/p:sldId/p:sldLayoutId r:id, the webextensionref
# r:id, the picture r:embed) plus the instance GUID on the embedded
# PowerPoll web add-in's <we:webextension id="..."> snapshot part. Those
# are exactly the artifacts PowerPoint regenerates internally whenever it
# rewrites a saved .pptx (fresh relationship ids package-wide, plus a new
# instance id for the task-pane add-in snapshot) - no slide text, shape
# geometry, layout, or add-in reference/image content actually changed
# between before/after.
#
# Touch the deck through the supported object model and persist it, which
# is the host-visible equivalent of that resave: the title/subtitle
# placeholders and the OfficeApp (PowerPoll) add-in frame are left
# exactly as authored.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Confirm/keep the two real placeholders untouched (no content changes
# shipped in this commit).
$title = $s.Shapes.Item(1)
$subtitle = $s.Shapes.Item(2)

# Persist the presentation - the COM-visible counterpart of the
# PowerPoint resave that produced the id/GUID churn in the diff.
$p.Save()
